$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unstyled reference cell (no explicit NumberFormat/border) used to restore the
# default style on cells where we temporarily force Text formatting so that
# numeric-looking values (e.g. "1.00", "0.999") stay literal text, matching the
# original inline-string cells instead of being coerced into numbers.
$plainStyle = $ws.Range("C2").Style

$ws.Range("D2").Value = '61.926.04'
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").Value = '3.413.96'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.59'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.481'
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = '  +0.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.04'
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = '  +4.49%  '
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("E11").Value = '  +1.74%  '
$ws.Range("D12").Value = '3.998.81'
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.29'
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = '  -4.75%  '
$ws.Range("D15").Value = '3.412.90'
$ws.Range("E15").Value = '  -1.34%  '
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").Value = '61.909.68'
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.35'
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = '  +0.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.39'
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.85'
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = '  -4.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.04'
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = '  -2.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.564'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.07'
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '3.554.93'
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("E26").Value = '  -3.83%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.60'
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("E30").Value = '  -3.87%  '
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.33'
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = '  -2.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.02'
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = '  -2.51%  '
$ws.Range("E35").Value = '  +2.69%  '
$ws.Range("E36").Value = '  +1.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '169.12'
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("E38").Value = '  -3.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '30.84'
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = '  -4.01%  '
$ws.Range("D40").Value = '3.446.59'
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("E41").Value = '  +1.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.58'
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("E45").Value = '  -3.26%  '
$ws.Range("E46").Value = '  -5.38%  '
$ws.Range("D47").Value = '2.537.05'
$ws.Range("E47").Value = '  -3.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.86'
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = '  +1.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.34'
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  -2.80%  '
$ws.Range("E51").Value = '  -5.75%  '
